# Update cryptos list values (prices and 1h volume %) pulled from coinranking.com
# Commit: Updated cryptos list on Wed Jun  5 04:42:30 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-numeric-looking price strings to stay text (matches source inlineStr cells)
$ws.Range("D2").Value = '71.190.33'
$ws.Range("E2").Value = '  +2.84%  '
$ws.Range("D3").Value = '3.810.50'
$ws.Range("E3").Value = '  +0.96%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '698.66'
$ws.Range("E5").Value = '  +10.57%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '173.05'
$ws.Range("E6").Value = '  +3.67%  '
$ws.Range("D7").Value = '3.810.65'
$ws.Range("E7").Value = '  +0.98%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.526'
$ws.Range("E9").Value = '  +0.93%  '
$ws.Range("E10").Value = '  +3.00%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.52'
$ws.Range("E11").Value = '  +11.29%  '
$ws.Range("E12").Value = '  +0.61%  '
$ws.Range("E13").Value = '  +8.77%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.35'
$ws.Range("E14").Value = '  +3.51%  '
$ws.Range("D15").Value = '4.448.89'
$ws.Range("E15").Value = '  +0.85%  '
$ws.Range("D16").Value = '3.823.33'
$ws.Range("E16").Value = '  +1.45%  '
$ws.Range("D17").Value = '71.087.20'
$ws.Range("E17").Value = '  +2.67%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.89'
$ws.Range("E18").Value = '  +1.48%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.23'
$ws.Range("E19").Value = '  +2.94%  '
$ws.Range("E20").Value = '  +1.12%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.23'
$ws.Range("E21").Value = '  +17.86%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '481.46'
$ws.Range("E22").Value = '  +4.05%  '
$ws.Range("E23").Value = '  +1.46%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.02'
$ws.Range("E24").Value = '  +1.79%  '
$ws.Range("E25").Value = '  +0.62%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.40'
$ws.Range("E26").Value = '  +2.78%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.18'
$ws.Range("E27").Value = '  +1.38%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.49'
$ws.Range("E28").Value = '  +4.13%  '
$ws.Range("D29").Value = '3.958.45'
$ws.Range("E29").Value = '  +0.89%  '
$ws.Range("E30").Value = '  -0.04%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.08'
$ws.Range("E31").Value = '  +14.50%  '
$ws.Range("E32").Value = '  -0.25%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.56'
$ws.Range("E33").Value = '  +6.73%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '29.61'
$ws.Range("E34").Value = '  +3.95%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.176'
$ws.Range("E35").Value = '  -0.65%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.24'
$ws.Range("E36").Value = '  +3.31%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  +0.06%  '
$ws.Range("D38").Value = '3.757.13'
$ws.Range("E38").Value = '  +0.82%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.103'
$ws.Range("E39").Value = '  +1.66%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.49'
$ws.Range("E40").Value = '  +5.74%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.99'
$ws.Range("E41").Value = '  +3.30%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.22'
$ws.Range("E42").Value = '  +11.23%  '
$ws.Range("B43").Value = 'Mantle'
$ws.Range("C43").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.968'
$ws.Range("E43").Value = '  +0.53%  '
$ws.Range("B44").Value = 'FLOKI'
$ws.Range("C44").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.000324'
$ws.Range("E44").Value = '  +21.78%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.999'
$ws.Range("E45").Value = '  -0.14%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '45.56'
$ws.Range("E47").Value = '  +4.85%  '
$ws.Range("B48").Value = 'OKB'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '49.32'
$ws.Range("E48").Value = '  +4.95%  '
$ws.Range("B49").Value = 'Monero'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '160.27'
$ws.Range("E49").Value = '  +1.31%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.42'
$ws.Range("E51").Value = '  +1.28%  '
